# Applies three related edits described by the commit:
#  1. Add the demo-video URL as a new run right after the "Demo Video URL\t: " line.
#  2. Append a short comment ". LiveData used when it is useful." (with the same
#     sz/szCs run formatting used throughout this bullet list) to the end of the
#     "MVVM architecture: ..." bullet, and move the "_GoBack" bookmark there
#     (that's where Word leaves it after the most recent edit).
#  3. Remove the old "_GoBack" bookmark from the "Git: ... spectacular." bullet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Demo Video URL: add the link run after the ": " run.
# ---------------------------------------------------------------------------
$videoPara = $d.Paragraphs.Item(2)
$videoRange = $videoPara.Range
$videoInsertPoint = $d.Range($videoRange.End - 1, $videoRange.End - 1)
$videoInsertPoint.InsertAfter("https://drive.google.com/open?id=1fZyECWsLUazJmvlSGGkXAdm-gpT8aocV")

# ---------------------------------------------------------------------------
# 2) MVVM architecture bullet: append the LiveData sentence + move _GoBack here.
#    We rebuild the paragraph's XML (same paraId/rsids/run rsids) so the only
#    observable change is the appended run + bookmark.
# ---------------------------------------------------------------------------
$mvvmPara = $d.Paragraphs.Item(70)
$mvvmRange = $mvvmPara.Range

$mvvmXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
  '<w:p w14:paraId="37E52DCA" w14:textId="087E482E" w:rsidR="00990775" w:rsidRPr="00B15944" w:rsidRDefault="0084049A" w:rsidP="0084049A">' +
    '<w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00B15944"><w:rPr><w:b/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>MVVM architecture</w:t></w:r>' +
    '<w:r w:rsidRPr="00B15944"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>: data stored in viewModel, so rotating &amp; recreating views do not cause errors</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>. LiveData used when it is useful.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$mvvmRange.InsertXML($mvvmXml)

# ---------------------------------------------------------------------------
# 3) Remove the stale _GoBack bookmark from the "Git: ... spectacular." bullet.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
